$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows above row 2 (before the 2025 row), shifting existing
# data down.
$insertRange = $ws.Range("A2:A6").EntireRow
$insertRange.Insert()

# After the insert, the old row 2 (with its formatting) has moved down to
# row 7. Re-apply that row's formatting to the 5 freshly inserted rows.
$formatSource = $ws.Range("A7:I7")
$formatTarget = $ws.Range("A2:I6")
$formatSource.Copy()
$formatTarget.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new rows for years 2020-2024 with the "Roth IRA ctrb" (column G)
# values, leaving the other columns blank.
$years = @(2020, 2021, 2022, 2023, 2024)
$rothIraCtrb = @(7000, 7000, 7500, 7500, 7500)

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 7).Value = $rothIraCtrb[$i]
}

# Update the active selection to match the final state (cell G7 selected).
$ws.Range("G7").Select()
